# INTERNSHIP MODULES LIST UPDATED
# Append three newly-added internship modules to the bottom of the table
# (columns: A=Code, B=Name, C=Credits, D=Graded).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 30 - MEEN40760 / ME MSE PWE (Long)
$ws.Range("B30").Value = "ME MSE PWE (Long)"
$ws.Range("A30").Value = "MEEN40760"
$ws.Range("C30").Value = 30
$ws.Range("D30").Value = "N"

# Row 31 - MEEN40930 / Professional Work Placement
$ws.Range("B31").Value = "Professional Work Placement"
$ws.Range("A31").Value = "MEEN40930"
$ws.Range("C31").Value = 20
$ws.Range("D31").Value = "N"

# Row 32 - BSEN40230 / ME Professional Work Experience (name reuses an
# existing shared string, so write it after the brand-new A32 code).
$ws.Range("A32").Value = "BSEN40230"
$ws.Range("B32").Value = "ME Professional Work Experience"
$ws.Range("C32").Value = 30
$ws.Range("D32").Value = "N"

# Leave the sheet's selection/scroll position where the author left it
# after entering the new rows.
$ws.Range("F31").Select()
$excel.ActiveWindow.ScrollRow = 14
